$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on cells whose new values look numeric (e.g. "7.13"),
# so Excel keeps them as text strings instead of converting to a Number,
# matching the original inlineStr (text) cell type from the source diff.
$textCells = @("D5","D6","D7","D9","D10","D11","D13","D15","D17","D19","D21","D22","D23","D27","D28","D29","D30","D31","D32","D34","D35","D36","D40","D41","D43","D44","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + volume%), row by row.

# Row 2
$ws.Range("D2").Value = "43.055.02"
$ws.Range("E2").Value = "  +2.65%  "

# Row 3
$ws.Range("D3").Value = "2.302.02"
$ws.Range("E3").Value = "  +1.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "310.23"
$ws.Range("E5").Value = "  +1.70%  "

# Row 6
$ws.Range("D6").Value = "100.50"
$ws.Range("E6").Value = "  +5.46%  "

# Row 7
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  +1.84%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  +6.25%  "

# Row 10
$ws.Range("D10").Value = "36.31"
$ws.Range("E10").Value = "  +3.64%  "

# Row 11
$ws.Range("D11").Value = "0.0823"
$ws.Range("E11").Value = "  +4.38%  "

# Row 12
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("D13").Value = "7.13"
$ws.Range("E13").Value = "  +7.64%  "

# Row 14
$ws.Range("D14").Value = "2.659.54"
$ws.Range("E14").Value = "  +1.81%  "

# Row 15
$ws.Range("D15").Value = "14.96"
$ws.Range("E15").Value = "  +3.88%  "

# Row 16
$ws.Range("D16").Value = "2.301.55"
$ws.Range("E16").Value = "  +2.19%  "

# Row 17
$ws.Range("D17").Value = "0.807"
$ws.Range("E17").Value = "  +2.02%  "

# Row 18
$ws.Range("D18").Value = "42.990.99"
$ws.Range("E18").Value = "  +2.71%  "

# Row 19
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  +2.88%  "

# Row 21
$ws.Range("D21").Value = "6.07"
$ws.Range("E21").Value = "  +1.94%  "

# Row 22
$ws.Range("D22").Value = "68.22"
$ws.Range("E22").Value = "  +0.82%  "

# Row 23
$ws.Range("D23").Value = "240.09"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24
$ws.Range("E24").Value = "  +4.91%  "

# Row 25
$ws.Range("E25").Value = "  +1.63%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("D27").Value = "24.37"
$ws.Range("E27").Value = "  +3.00%  "

# Row 28
$ws.Range("D28").Value = "38.77"
$ws.Range("E28").Value = "  +6.15%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  +2.39%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.66"
$ws.Range("E30").Value = "  +1.90%  "

# Row 31
$ws.Range("D31").Value = "168.35"
$ws.Range("E31").Value = "  +5.17%  "

# Row 32
$ws.Range("D32").Value = "5.34"
$ws.Range("E32").Value = "  +2.57%  "

# Row 33
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  -1.24%  "

# Row 35
$ws.Range("D35").Value = "17.82"
$ws.Range("E35").Value = "  +4.37%  "

# Row 36
$ws.Range("D36").Value = "0.0740"
$ws.Range("E36").Value = "  +0.53%  "

# Row 37
$ws.Range("E37").Value = "  +0.29%  "

# Row 38
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("E39").Value = "  +2.06%  "

# Row 40
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +1.09%  "

# Row 41
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  +5.97%  "

# Row 42
$ws.Range("E42").Value = "  -4.97%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0289"
$ws.Range("E43").Value = "  +2.30%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.28"
$ws.Range("E44").Value = "  +2.39%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.967.74"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +3.07%  "

# Row 47
$ws.Range("D47").Value = "9.83"
$ws.Range("E47").Value = "  -0.50%  "

# Row 48
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  +16.72%  "

# Row 49
$ws.Range("D49").Value = "55.29"
$ws.Range("E49").Value = "  +4.10%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.55"
$ws.Range("E50").Value = "  +2.69%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.526.00"
$ws.Range("E51").Value = "  +1.65%  "
